$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B-G, rows 2-10 (regenerated s_val data filtering save games)
$data = @{
    2  = @(0.1554434735375247, 0.3375848360084654, 16.98373111632243, 0.4998867070740569, 1, 17.97664613294248)
    3  = @(3.182878228561681, 86.29678392075563, 0.1529057820181812, 6.48142807727062, 1, 96.11399600860611)
    4  = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 6.48142807727062, 0, 9.793184359356808)
    5  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729)
    6  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    7  = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 6.741336633845642)
    8  = @(3.182878228561681, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 0, 4.733082622659194)
    9  = @(0.06328177979961902, 1.65323645889881, 0.1529057820181812, 6.48142807727062, 0, 8.35085209798723)
    10 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
